# ---------------------------------------------------------------------------
# Converts an "RRGGBB" hex string into the packed BGR integer that the
# PowerPoint object model expects for ColorFormat/ThemeColor .RGB values.
# ---------------------------------------------------------------------------
function HexToBGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 5 contains a table ("Google Shape;122;p17"); point it at a
#    different built-in table style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{2038C22E-64FF-4DB0-AC6F-C0D54EB73065}")

# ---------------------------------------------------------------------------
# 2. Re-colour the deck's theme so its colour scheme matches the standard
#    Office palette instead of the old Integral "Red Violet" palette.
#    All slides share the one slide master/theme, so touching it via any
#    single slide updates the theme for the whole deck.
# ---------------------------------------------------------------------------
$cs = $p.Slides.Item(1).ThemeColorScheme

$cs.Item(1).RGB  = HexToBGR "000000"   # dk1
$cs.Item(2).RGB  = HexToBGR "FFFFFF"   # lt1
$cs.Item(3).RGB  = HexToBGR "44546A"   # dk2
$cs.Item(4).RGB  = HexToBGR "E7E6E6"   # lt2
$cs.Item(5).RGB  = HexToBGR "5B9BD5"   # accent1
$cs.Item(6).RGB  = HexToBGR "ED7D31"   # accent2
$cs.Item(7).RGB  = HexToBGR "A5A5A5"   # accent3
$cs.Item(8).RGB  = HexToBGR "FFC000"   # accent4
$cs.Item(9).RGB  = HexToBGR "4472C4"   # accent5
$cs.Item(10).RGB = HexToBGR "70AD47"   # accent6
$cs.Item(11).RGB = HexToBGR "0563C1"   # hlink
$cs.Item(12).RGB = HexToBGR "954F72"   # folHlink
